$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 10 (pushes "fossil_routes" and everything below it
# down by one row), then populate it with the new
# "chemical_recycling_pyrolysis" parameter, right after
# "chemical_recycling_gasification" (row 9).
$ws.Rows("10:10").Insert()

$ws.Range("A10").Value = "chemical_recycling_pyrolysis"
$ws.Range("B10").Value = $true
